$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ferroviario figure for 2024 Q4 is now available (was "ND"); set the real value.
$ws.Range("E5").Value = 14238.468000000001

# Update the "last updated" note.
$ws.Range("B33").Value = "Actualización: Enero 2025."

# The "ND No Disponible" legend row is no longer needed; remove it and shift
# the remaining notes/source rows up.
$ws.Rows(34).Delete()
